$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 341, shifting the existing
# rows 341-428 down to 343-430.
$ws.Range("A341:A342").EntireRow.Insert()

# New row 341: Murcott / Primera, Región de O'Higgins
$ws.Range("A341").Value = 10
$ws.Range("B341").Value = "Vega Modelo de Temuco"
$ws.Range("C341").Value = "La Araucanía"
$ws.Range("D341").Value = 44508
$ws.Range("E341").Value = 9
$ws.Range("F341").Value = "Fruta"
$ws.Range("G341").Value = 100102
$ws.Range("H341").Value = "Cítricos"
$ws.Range("I341").Value = 100102004
$ws.Range("J341").Value = "Mandarina"
$ws.Range("K341").Value = "Murcott"
$ws.Range("L341").Value = "Primera"
$ws.Range("M341").Value = 250
$ws.Range("N341").Value = 10000
$ws.Range("O341").Value = 10000
$ws.Range("P341").Value = 10000
$ws.Range("Q341").Value = "$/caja 10 kilos"
$ws.Range("R341").Value = "Región de O'Higgins"
$ws.Range("S341").Value = 1000
$ws.Range("T341").Value = 10

# New row 342: Murcott / Tercera, Región de O'Higgins
$ws.Range("A342").Value = 10
$ws.Range("B342").Value = "Vega Modelo de Temuco"
$ws.Range("C342").Value = "La Araucanía"
$ws.Range("D342").Value = 44508
$ws.Range("E342").Value = 9
$ws.Range("F342").Value = "Fruta"
$ws.Range("G342").Value = 100102
$ws.Range("H342").Value = "Cítricos"
$ws.Range("I342").Value = 100102004
$ws.Range("J342").Value = "Mandarina"
$ws.Range("K342").Value = "Murcott"
$ws.Range("L342").Value = "Tercera"
$ws.Range("M342").Value = 125
$ws.Range("N342").Value = 5400
$ws.Range("O342").Value = 5400
$ws.Range("P342").Value = 5400
$ws.Range("Q342").Value = "$/bandeja 18 kilos"
$ws.Range("R342").Value = "Región de O'Higgins"
$ws.Range("S342").Value = 300
$ws.Range("T342").Value = 18
